$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New question/answer data for rows 2-8 (A: Question, B: User Answer, C: Correct Answer)
$data = @(
    @("What is the standard algorithm for inserting an item into a sorted array?", "Insertion sort", "Insertion sort"),
    @("What is a sequential file?", "Sequential file is a file that stores data in a sequential manner.", "Sequential file is a file that stores data in a sequential manner."),
    @("What is the standard algorithm for sorting an array?", "Selection sort", "Selection sort"),
    @("What is the standard algorithm for deleting an item from a sorted array?", "Deletion sort", "Deletion sort"),
    @("What is the standard algorithm for appending to a sequential file?", "OPEN and WRITE", "OPEN and WRITE"),
    @("What is the standard algorithm for searching an unsorted array?", "Binary search", "Linear search"),
    @("What is the difference between a sequential file and a relative file?", "Sequential files can only be read from beginning to end, while relative files can be read from any point.", "Sequential files can only be read from beginning to end, while relative files can be read from any point.")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove column D entirely
$ws.Range("D1:D11").Delete()

# Remove now-unused rows 9, 10, 11 (below row 8)
$ws.Range("A9:C11").Delete()
